$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Base job description paragraph used throughout this workbook.
$para = "We are seeking a Software Engineer to build and maintain high-quality software solutions.`nWork with global teams to drive innovation and deliver scalable applications.`nJoin Akkodis and be part of a tech-driven, collaborative environment."

# The new job posting's description is this paragraph repeated 5 times.
$desc = $para + $para + $para + $para + $para

# Append a new row (row 10) with the new job posting JD_009.
$ws.Range("A10").Value = "JD_009"
$ws.Range("B10").Value = "Senior Java Engineer"
$ws.Range("C10").Value = $desc
$ws.Range("D10").Value = 1
$ws.Range("E10").Value = 5

# Avoid leaving a stale explicit row height from the long, multi-line
# description text - keep the new row using the sheet's default height,
# consistent with every other data row.
$ws.Rows.Item(10).EntireRow.AutoFit()
